$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (Oct 06, 2022 NYSE/NASDAQ market diary) to be duplicated into rows 5 and 6
$rowValues = @(
    "Oct 06, 2022",
    "NYSE",
    3346,
    974,
    2242,
    130,
    32,
    122,
    336249580,
    576118947,
    919538564,
    0.7,
    5396,
    1635873706,
    2625683225,
    "NASDAQ",
    4293189892,
    4878,
    1713,
    2902,
    263,
    59,
    199,
    0.9399999999999999,
    19974,
    1460424411,
    2334144952,
    4091868336
)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(5, $col).Value = $rowValues[$i]
    $ws.Cells.Item(6, $col).Value = $rowValues[$i]
}
